# REVER_DailyTracker_BALAJI.xlsx - "Add files via upload" edit
# Fills in Application / Task details for several rows of the APR-2021 daily
# tracker sheet, adding the corresponding shared strings, and moves the
# active selection to D14 to match the saved state of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content cells -----------------------------------------------------
# Cells are written in the same order the new shared strings were appended
# to xl/sharedStrings.xml so the unique-string table comes out identical.

# Row 9 (No. 8)
$ws.Range("C9").Value = "NMVar"
$ws.Range("D9").Value = "Issue fixed"

# Row 8 (No. 7)
$ws.Range("D8").Value = "Again some changes given by Vijay san so it is going on"

# Row 10 (No. 9)
$ws.Range("C10").Value = "Leave"
$ws.Range("D10").Value = "Personal Leave"

# Row 16 (No. 15)
$ws.Range("D16").Value = "Point number 20 going on"

# Row 13 (No. 12)
$ws.Range("D13").Value = "Messages taken from all screens point number 12 fixed"

# Row 15 (No. 14)
$ws.Range("D15").Value = "Tamil New Year"

# Row 14 (No. 13)
$ws.Range("D14").Value = "Telugu New Year"

# --- Cells re-using already existing shared strings -------------------------

# Row 7 (No. 6)
$ws.Range("C7").Value = "Mujistore "
$ws.Range("D7").Value = "Fixing going onMujistore submit button issue in video call menu"
# This cell already existed (empty) with a style lacking the thin border
# formatting applied to the rest of the column; restore the border so the
# cell matches its neighbours.
$ws.Range("D7").Borders.LineStyle = 1

# Row 8 (No. 7)
$ws.Range("C8").Value = "Mujistore "

# Row 13 (No. 12)
$ws.Range("C13").Value = "Mujistore "

# Row 14 (No. 13)
$ws.Range("C14").Value = "Holiday"

# Row 15 (No. 14)
$ws.Range("C15").Value = "Holiday"

# Row 16 (No. 15)
$ws.Range("C16").Value = "NMVar"

# --- Selection ---------------------------------------------------------------
$ws.Range("D14").Select() | Out-Null
